# Generate Report for Handback
# The handback file (9014c13d-86bc-497e-8618-682a271afbd8.md) has been
# successfully handed back and is now in sync with en-US. Update the
# status, handback datetime, and clear the stale error detail across
# the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet -------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# --- zh-cn sheet ------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $newStatus
$zhcn.Range("K3").Value = "2016-08-18 12:48:50"
$zhcn.Range("P3").Value = ""

# --- de-de sheet ------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $newStatus
$dede.Range("K3").Value = "2016-08-18 12:48:57"
$dede.Range("P3").Value = ""
